$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 157, shifting the existing rows (157-226) down to (158-227)
$ws.Rows("157:157").Insert()

# Populate the new row 157 with the new record.
# Columns A, B, C, E, F, G, H, I, N, Q, R are identical to the rest of the
# "Macroferia Regional de Talca - Zapallo italiano" records in this block.
$ws.Range("A157").Value = 5
$ws.Range("B157").Value = "Macroferia Regional de Talca"
$ws.Range("C157").Value = "Maule"
$ws.Range("D157").Value = 44510
$ws.Range("E157").Value = 7
$ws.Range("F157").Value = 100112032
$ws.Range("G157").Value = "Zapallo italiano"
$ws.Range("H157").Value = "Sin especificar"
$ws.Range("I157").Value = "Primera"
$ws.Range("J157").Value = 500
$ws.Range("K157").Value = 7000
$ws.Range("L157").Value = 7000
$ws.Range("M157").Value = 7000
$ws.Range("N157").Value = "$/caja 60 unidades"
$ws.Range("O157").Value = "Región del Maule"
$ws.Range("P157").Value = 117
$ws.Range("Q157").Value = 60
$ws.Range("R157").Value = "Hortaliza"
